$wb = $excel.ActiveWorkbook

# Add the new worksheet at the end of the workbook and name it "cheese"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "cheese"

$ws = $wb.Worksheets.Item("cheese")

# Header row
$ws.Range("A1").Value = "COMP 222  B"
$ws.Range("B1").Value = "DATA STRUC/ALGO"
$ws.Range("C1").Value = "INTRO TO DATA STRUCT & ALGORITHMS"
$ws.Range("D1").Value = "11:30:00"
$ws.Range("E1").Value = "12:45:00"
$ws.Range("F1").Value = "TR"
$ws.Range("G1").Value = "STEM"

# Room number looks numeric, so force it to Text first (otherwise Excel
# stores it as a number instead of a string), then clear the resulting Text
# number-format back off the cell once the string value has been committed.
$ws.Range("H1").NumberFormat = "@"
$ws.Range("H1").Value = "376"
$ws.Range("H1").ClearFormats()

# Second row
$ws.Range("A2").Value = "EDUC 312  A"
$ws.Range("B2").Value = "MUSIC METHODS"
$ws.Range("C2").Value = "SECONDARY MUSIC METH-MUS ED MAJ"
$ws.Range("D2").Value = "10:05:00"
$ws.Range("E2").Value = "11:20:00"
$ws.Range("F2").Value = "TR"
$ws.Range("G2").Value = "PFAC"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "68"
$ws.Range("H2").ClearFormats()

$wb.Save()
